$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(41, "14265594", "2025-08-05", "Alex Michelsen", "Karen Khachanov", "Gana Alex Michelsen", 2.75),
    @(42, "14266317", "2025-08-04", "Marta Kostyuk", "Elena Rybakina", "Gana Marta Kostyuk", 3.75),
    @(43, "14266318", "2025-08-04", "Victoria Mboko", "Jessica Bouzas Maneiro", "Gana Jessica Bouzas Maneiro", 3.5),
    @(44, "14339494", "2025-08-05", "Gonzalo Bueno", "Maxim Mrva", "Gana Gonzalo Bueno", 2.2),
    @(45, "14340604", "2025-08-04", "Naoki Nakagawa", "Alfredo Perez", "Gana Naoki Nakagawa", 2.63),
    @(46, "14340605", "2025-08-04", "Yibing Wu", "Andre Ilagan", "Gana Andre Ilagan", 5.5),
    @(47, "14344348", "2025-08-04", "Alex Molcan", "Martin Krumich", "Gana Martin Krumich", 3.4),
    @(48, "14266317", "2025-08-04", "Marta Kostyuk", "Elena Rybakina", "Gana Marta Kostyuk", 3.75)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = ""
    $ws.Cells.Item($r, 8).Value = ""
}
